$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks -------------------------------------------------
# B1 was "ACTION" (shared string index 1); it becomes "TERMINALS" (new
# shared string). Preserve the old "ACTION" string by also writing it into
# the new I1 header cell below, so it is not dropped from the table.
$ws.Range("B1").Value = "TERMINALS"

# New "ACTION" mini-header over the new TERMINALS column (I), merged I1:I2,
# matching the look of the existing G1:H1 "GOTO" header (style copied from H1).
# Merge first, then apply formats, so Excel doesn't split the thin box
# border into separate top/bottom pieces for the merged range.
$ws.Range("I1:I2").Merge() | Out-Null

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "ACTION"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New empty bold/vertically-centered cell to the right of the new column.
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").VerticalAlignment = -4108        # xlCenter

# --- Row 7: fill in previously-empty cells ------------------------------
# (Set D7 first so the new "R3" shared string is introduced ahead of "N/A".)
$ws.Range("D7").Value = "R3"
$ws.Range("E7").Value = "R3"
$ws.Range("F7").Value = "R3"

# --- New "TERMINALS" data column (I), rows 3-8 --------------------------
$ws.Range("F3").Copy() | Out-Null
$ws.Range("I3:I6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I3").Value = "N/A"
$ws.Range("I4").Value = "N/A"
$ws.Range("I5").Value = "N/A"
$ws.Range("I6").Value = "N/A"

# --- Row 7 (continued) ---------------------------------------------------
$ws.Range("B7").Value = "ERR"
$ws.Range("C7").Value = "ERR"
$ws.Range("G7").Value = "ERR"
$ws.Range("H7").Value = "ERR"
$ws.Range("I7").Value = "R3"

# --- Row 8: new state row (I5), with its own data -----------------------
$ws.Range("A8").Value = "I5"
$ws.Range("B8").Value = "ERR"
$ws.Range("C8").Value = "ERR"
$ws.Range("D8").Value = "R4"
$ws.Range("E8").Value = "R4"
$ws.Range("F8").Value = "R4"
$ws.Range("G8").Value = "ERR"
$ws.Range("H8").Value = "ERR"
$ws.Range("I8").Value = "R4"

# --- Selection / active cell --------------------------------------------
$ws.Range("A9").Select() | Out-Null
